$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.173.30'
$ws.Range("E2").Value = '  +3.37%  '
$ws.Range("D3").Value = '1.904.45'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'326.07"
$ws.Range("E5").Value = '  +3.34%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = "'0.5153"
$ws.Range("E7").Value = '  +0.90%  '
$ws.Range("D8").Value = "'0.3992"
$ws.Range("E8").Value = '  +1.47%  '
$ws.Range("D9").Value = "'0.08455"
$ws.Range("E9").Value = '  +0.51%  '
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").Value = "'23.23"
$ws.Range("E12").Value = '  +13.00%  '
$ws.Range("D13").Value = "'6.439"
$ws.Range("E13").Value = '  +2.97%  '
$ws.Range("D14").Value = '1.910.12'
$ws.Range("E14").Value = '  +1.46%  '
$ws.Range("D15").Value = "'7.349"
$ws.Range("E15").Value = '  +0.03%  '
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = "'94.89"
$ws.Range("E17").Value = '  +1.87%  '
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("D19").Value = "'0.06620"
$ws.Range("E19").Value = '  -1.35%  '
$ws.Range("D20").Value = "'18.37"
$ws.Range("E20").Value = '  +2.83%  '
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").Value = "'5.993"
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").Value = '30.182.87'
$ws.Range("E23").Value = '  +3.38%  '
$ws.Range("D24").Value = "'11.24"
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("D25").Value = "'2.210"
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("D26").Value = '2.135.75'
$ws.Range("E26").Value = '  +1.69%  '
$ws.Range("D27").Value = "'21.59"
$ws.Range("E27").Value = '  +3.08%  '
$ws.Range("D28").Value = "'162.15"
$ws.Range("E28").Value = '  +2.22%  '
$ws.Range("D29").Value = "'2.384"
$ws.Range("E29").Value = '  -2.06%  '
$ws.Range("D30").Value = "'129.60"
$ws.Range("E30").Value = '  +1.99%  '
$ws.Range("D31").Value = "'1.095"
$ws.Range("E31").Value = '  +3.24%  '
$ws.Range("E32").Value = '  +1.48%  '
$ws.Range("D33").Value = "'6.072"
$ws.Range("E33").Value = '  +3.30%  '
$ws.Range("D34").Value = "'3.656"
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("E35").Value = '  +0.73%  '
$ws.Range("D36").Value = "'0.06573"
$ws.Range("E36").Value = '  -0.36%  '
$ws.Range("D37").Value = "'0.2203"
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("D38").Value = "'5.184"
$ws.Range("E38").Value = '  +1.37%  '
$ws.Range("D39").Value = "'1.226"
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("E40").Value = '  +5.82%  '
$ws.Range("D41").Value = "'8.794"
$ws.Range("E41").Value = '  -3.07%  '
$ws.Range("D42").Value = "'0.6517"
$ws.Range("E42").Value = '  +1.09%  '
$ws.Range("D43").Value = "'1.233"
$ws.Range("E43").Value = '  -0.52%  '
$ws.Range("D44").Value = "'0.6130"
$ws.Range("E44").Value = '  +1.43%  '
$ws.Range("D45").Value = "'13.16"
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("D46").Value = "'3.720"
$ws.Range("E47").Value = '  +1.25%  '
$ws.Range("D48").Value = "'1.244"
$ws.Range("E48").Value = '  +0.99%  '
$ws.Range("D49").Value = "'125.01"
$ws.Range("E49").Value = '  +1.67%  '
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("D51").Value = "'79.22"
$ws.Range("E51").Value = '  +2.37%  '
